# Update "想去人数" (want-to-go count) figures on the "展览" and "全部类型"
# sheets to reflect newly generated output.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 13038
$ws1.Range("F4").Value = 315
$ws1.Range("F5").Value = 638
$ws1.Range("F6").Value = 207
$ws1.Range("F7").Value = 407
$ws1.Range("F8").Value = 1233

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 13038
$ws4.Range("F5").Value = 315
$ws4.Range("F6").Value = 638
$ws4.Range("F7").Value = 207
$ws4.Range("F10").Value = 407
$ws4.Range("F11").Value = 1233
